$wb = $excel.ActiveWorkbook
$wsElim = $wb.Worksheets.Item("HypothyroidismEliminate")
$wsElim.Range("J2").Value = "[Diabetic, PCOS, Diabetic, PCOS]"
$wsElim.Range("D3").Value = "Veg Veg "
$wsElim.Range("D4").Value = "Veg Veg "
$wsElim.Range("D5").Value = "Veg Veg "
$wsElim.Range("D6").Value = "Veg Veg "
$wsElim.Range("J6").Value = "[PCOS, PCOS]"
$wsElim.Range("D8").Value = "Veg Veg "
$wsElim.Range("J8").Value = "[Diabetic, Diabetic]"
$wsElim.Range("C10").Value = "[Snacks, Breakfast, Snacks, Breakfast]"
$wsElim.Range("D11").Value = "Veg Veg "
$wsElim.Range("C12").Value = "[Dinner, Lunch, Dinner, Lunch]"
$wsElim.Range("D12").Value = "Veg Veg "
$wsElim.Range("D13").Value = "Veg Veg "
$wsElim.Range("D14").Value = "Veg Veg "
$wsElim.Range("C15").Value = "[Breakfast, Breakfast]"
$wsElim.Range("D17").Value = "Vegetarian Veg Vegetarian Veg "
$wsElim.Range("J18").Value = "[Diabetic, Diabetic]"
$wsElim.Range("D19").Value = "Veg Veg "
$wsElim.Range("J19").Value = "[Diabetic, Diabetic]"
$wsElim.Range("C20").Value = "[Snacks, Snacks]"
$wsElim.Range("D20").Value = "Jain Jain "
$wsElim.Range("C21").Value = "[Snacks, Snacks]"
$wsElim.Range("D21").Value = "Veg Veg "
$wsElim.Range("C22").Value = "[Breakfast, Breakfast]"
$wsElim.Range("D22").Value = "Veg Veg "
$wsElim.Range("C23").Value = "[Breakfast, Breakfast]"
$wsElim.Range("J24").Value = "[PCOS, PCOS]"
$wsElim.Range("D25").Value = "Veg Veg "
$wsElim.Range("J27").Value = "[Diabetic, Diabetic]"
$wsElim.Range("C30").Value = "[Breakfast, Breakfast]"
$wsElim.Range("D31").Value = "Veg Veg "
$wsElim.Range("C32").Value = "[Snacks, Snacks]"
$wsElim.Range("D32").Value = "Veg Veg "
$wsElim.Range("D33").Value = "Veg Veg "
$wsElim.Range("D34").Value = "Jain Jain "
$wsElim.Range("C35").Value = "[Snacks, Snacks]"
$wsElim.Range("D35").Value = "Veg Veg "
$wsElim.Range("C36").Value = "[Snacks, Breakfast, Snacks, Breakfast]"
$wsElim.Range("C37").Value = "[Snacks, Snacks]"
$wsElim.Range("C38").Value = "[Snacks, Snacks]"
$wsElim.Range("C39").Value = "[Snacks, Snacks]"
$wsElim.Range("C40").Value = "[Breakfast, Breakfast]"
$wsElim.Range("C41").Value = "[Dinner, Dinner]"
$wsElim.Range("D43").Value = "Veg Veg "

$wsElimination = $wb.Worksheets.Item("HypothyroidismElimination")
$wsElimination.Range("D2").Value = "[Veg, Veg, Veg]"
